$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (SOS_ON_DISPLAY_MANUFACTURER): expand manufacturer list, drop category/SPIRITS columns ---
$ws.Range("D4").Value = "PERNOD RICARD,DIAGEO,BEAM SUNTORY,PROXIMO,BACARDI,E&J GALLO,BROWN-FORMAN,CONSTELLATION,SAZERAC,HEAVEN HILL"
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()

# --- Row 5 (SOS_ON_DISPLAY_BRANDS): expand brand list ---
$ws.Range("D5").Value = "JAMESON,ABSOLUT,SEAGRAM'S,MALIBU,KAHLUA,GLENLIVET,CHIVAS REGAL,BEEFEATER,ALTOS,AVION,MARTELL,SMIRNOFF,TITO'S,NEW AMSTERDAM,CAPTAIN MORGAN,BACARDI,SAILOR JERRY,JACK DANIEL'S,JIM BEAM,MAKER'S MARK,JOSE CUERVO,SAUZA,1800,BAILEYS,JAGERMEISTER,CROWN ROYAL,FIREBALL,CANADIAN CLUB,BOMBAY,TANQUERAY,HENNESSY,COURVOISIER,REMY MARTIN,JOHNNIE WALKER"
$ws.Rows.Item(5).RowHeight = 102.2

# --- Row 6 (SHARE_OF_DISPLAY_MANUFACTURER): expand manufacturer list, drop category/SPIRITS columns ---
$ws.Range("D6").Value = "PERNOD RICARD,DIAGEO,BEAM SUNTORY,PROXIMO,BACARDI,E&J GALLO,BROWN-FORMAN,CONSTELLATION,SAZERAC,HEAVEN HILL"
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()

# --- Row 7 (SHARE_OF_DISPLAY_BRANDS): expand brand list ---
$ws.Range("D7").Value = "JAMESON,ABSOLUT,SEAGRAM'S,MALIBU,KAHLUA,GLENLIVET,CHIVAS REGAL,BEEFEATER,ALTOS,AVION,MARTELL,SMIRNOFF,TITO'S,NEW AMSTERDAM,CAPTAIN MORGAN,BACARDI,SAILOR JERRY,JACK DANIEL'S,JIM BEAM,MAKER'S MARK,JOSE CUERVO,SAUZA,1800,BAILEYS,JAGERMEISTER,CROWN ROYAL,FIREBALL,CANADIAN CLUB,BOMBAY,TANQUERAY,HENNESSY,COURVOISIER,REMY MARTIN,JOHNNIE WALKER"
$ws.Rows.Item(7).RowHeight = 102.2

# --- Row 8 (SOLO_SHARED): expand brand list ---
$ws.Range("D8").Value = "JAMESON,ABSOLUT,SEAGRAM'S,MALIBU,KAHLUA,GLENLIVET,CHIVAS REGAL,BEEFEATER,ALTOS,AVION,MARTELL,SMIRNOFF,TITO'S,NEW AMSTERDAM,CAPTAIN MORGAN,BACARDI,SAILOR JERRY,JACK DANIEL'S,JIM BEAM,MAKER'S MARK,JOSE CUERVO,SAUZA,1800,BAILEYS,JAGERMEISTER,CROWN ROYAL,FIREBALL,CANADIAN CLUB,BOMBAY,TANQUERAY,HENNESSY,COURVOISIER,REMY MARTIN,JOHNNIE WALKER"
$ws.Rows.Item(8).RowHeight = 102.2

# --- Update the active selection to E3 (template update) ---
$ws.Range("E3").Select()
